$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 27
$ws.Range('B27').Value = 6504313
$ws.Range('F27').Value = 'Guabira'
$ws.Range('G27').Value = 'Atletico Palmaflor Vinto'
$ws.Range('I27').Value = 0
$ws.Range('J27').Value = 'H'
$ws.Range('K27').Value = 1.75
$ws.Range('M27').Value = 4
$ws.Range('N27').Value = 2
$ws.Range('O27').Value = 3.6
$ws.Range('P27').Value = 3.6
$ws.Range('Q27').Value = -0.25
$ws.Range('R27').Value = 1.75
$ws.Range('S27').Value = 2.05
$ws.Range('T27').Value = 2.5
$ws.Range('U27').Value = 1.85
$ws.Range('V27').Value = 1.95
$ws.Range('W27').Value = 1
$ws.Range('Y27').Value = -1
$ws.Range('Z27').Value = 0.75
$ws.Range('AA27').Value = -1
$ws.Range('AB27').Value = -1
$ws.Range('AC27').Value = 0.95

# Row 28
$ws.Range('B28').Value = 6504831
$ws.Range('F28').Value = 'Libertad Gran Mamore FC'
$ws.Range('G28').Value = 'Always Ready'
$ws.Range('I28').Value = 2
$ws.Range('J28').Value = 'A'
$ws.Range('K28').Value = 3.3
$ws.Range('M28').Value = 1.909
$ws.Range('N28').Value = 3.1
$ws.Range('O28').Value = 3.5
$ws.Range('P28').Value = 2.2
$ws.Range('Q28').Value = 0.25
$ws.Range('R28').Value = 1.9
$ws.Range('S28').Value = 1.9
$ws.Range('T28').Value = 2.75
$ws.Range('U28').Value = 2.025
$ws.Range('V28').Value = 1.775
$ws.Range('W28').Value = -1
$ws.Range('Y28').Value = 1.2
$ws.Range('Z28').Value = -1
$ws.Range('AA28').Value = 0.8999999999999999
$ws.Range('AB28').Value = 0.5125
$ws.Range('AC28').Value = -0.5

# Row 142
$ws.Range('B142').Value = 7532413
$ws.Range('F142').Value = 'Libertad Gran Mamore FC'
$ws.Range('G142').Value = 'Club Aurora'
$ws.Range('H142').Value = 0
$ws.Range('J142').Value = 'A'
$ws.Range('K142').Value = 2.25
$ws.Range('L142').Value = 3.3
$ws.Range('M142').Value = 2.8
$ws.Range('N142').Value = 2.375
$ws.Range('O142').Value = 3.4
$ws.Range('P142').Value = 2.875
$ws.Range('Q142').Value = -0.25
$ws.Range('R142').Value = 2.025
$ws.Range('S142').Value = 1.775
$ws.Range('T142').Value = 2.5
$ws.Range('W142').Value = -1
$ws.Range('Y142').Value = 1.875
$ws.Range('Z142').Value = -1
$ws.Range('AA142').Value = 0.7749999999999999
$ws.Range('AB142').Value = -1
$ws.Range('AC142').Value = 0.8999999999999999

# Row 143
$ws.Range('B143').Value = 7532430
$ws.Range('F143').Value = 'Always Ready'
$ws.Range('G143').Value = 'Oriente Petrolero'
$ws.Range('H143').Value = 4
$ws.Range('J143').Value = 'H'
$ws.Range('K143').Value = 1.4
$ws.Range('L143').Value = 4.2
$ws.Range('M143').Value = 7
$ws.Range('N143').Value = 1.363
$ws.Range('O143').Value = 4.5
$ws.Range('P143').Value = 8.5
$ws.Range('Q143').Value = -1.5
$ws.Range('R143').Value = 2
$ws.Range('S143').Value = 1.8
$ws.Range('T143').Value = 3
$ws.Range('W143').Value = 0.363
$ws.Range('Y143').Value = -1
$ws.Range('Z143').Value = 1
$ws.Range('AA143').Value = -1
$ws.Range('AB143').Value = 0.8999999999999999
$ws.Range('AC143').Value = -1

# Row 148
$ws.Range('B148').Value = 7532419
$ws.Range('F148').Value = 'Oriente Petrolero'
$ws.Range('G148').Value = 'Jorge Wilstermann'
$ws.Range('H148').Value = 3
$ws.Range('K148').Value = 2.2
$ws.Range('L148').Value = 2.5
$ws.Range('M148').Value = 4.5
$ws.Range('N148').Value = 2.375
$ws.Range('O148').Value = 2.45
$ws.Range('P148').Value = 4.5
$ws.Range('Q148').Value = -0.25
$ws.Range('R148').Value = 1.9
$ws.Range('S148').Value = 1.9
$ws.Range('T148').Value = 2
$ws.Range('U148').Value = 1.95
$ws.Range('V148').Value = 1.85
$ws.Range('W148').Value = 1.375
$ws.Range('Z148').Value = 0.8999999999999999
$ws.Range('AB148').Value = 0.95
$ws.Range('AC148').Value = -1

# Row 150
$ws.Range('B150').Value = 7532421
$ws.Range('F150').Value = 'Guabira'
$ws.Range('G150').Value = 'Independiente Petrolero'
$ws.Range('H150').Value = 2
$ws.Range('K150').Value = 1.4
$ws.Range('L150').Value = 4.5
$ws.Range('M150').Value = 7.5
$ws.Range('N150').Value = 1.333
$ws.Range('O150').Value = 5.5
$ws.Range('P150').Value = 9.5
$ws.Range('Q150').Value = -1.5
$ws.Range('R150').Value = 1.85
$ws.Range('S150').Value = 1.95
$ws.Range('T150').Value = 3
$ws.Range('U150').Value = 1.825
$ws.Range('V150').Value = 1.975
$ws.Range('W150').Value = 0.333
$ws.Range('Z150').Value = 0.8500000000000001
$ws.Range('AB150').Value = -1
$ws.Range('AC150').Value = 0.9750000000000001

# Row 203
$ws.Range('B203').Value = 8010637
$ws.Range('E203').Value = 45381.66666666666
$ws.Range('F203').Value = 'San Antonio Bulo Bulo'
$ws.Range('G203').Value = 'Real Tomayapo'
$ws.Range('H203').Value = 3
$ws.Range('I203').Value = 0
$ws.Range('J203').Value = 'H'
$ws.Range('K203').Value = 1.727
$ws.Range('L203').Value = 3.5
$ws.Range('M203').Value = 4
$ws.Range('N203').Value = 1.666
$ws.Range('O203').Value = 3.8
$ws.Range('P203').Value = 5.25
$ws.Range('Q203').Value = -0.75
$ws.Range('R203').Value = 1.8
$ws.Range('S203').Value = 2
$ws.Range('T203').Value = 2.5
$ws.Range('U203').Value = 1.975
$ws.Range('V203').Value = 1.825
$ws.Range('W203').Value = 0.6659999999999999
$ws.Range('X203').Value = -1
$ws.Range('Y203').Value = -1
$ws.Range('Z203').Value = 0.8
$ws.Range('AA203').Value = -1
$ws.Range('AB203').Value = 0.9750000000000001
$ws.Range('AC203').Value = -1

# Row 204
$ws.Range('B204').Value = 8011587
$ws.Range('E204').Value = 45382.85416666666
$ws.Range('F204').Value = 'Club Aurora'
$ws.Range('G204').Value = 'Blooming'
$ws.Range('K204').Value = 1.533
$ws.Range('L204').Value = 3.75
$ws.Range('M204').Value = 5.5
$ws.Range('N204').Value = 1.55
$ws.Range('O204').Value = 3.8
$ws.Range('P204').Value = 6.5
$ws.Range('Q204').Value = -1
$ws.Range('U204').Value = 2
$ws.Range('V204').Value = 1.8

# Remove now-obsolete trailing rows (old rows 205 and 206)
$ws.Rows.Item(206).Delete()
$ws.Rows.Item(205).Delete()
